# Tests were added for annotation icon on delivery panel
#
# Appends a new "auAnnotationUser1" test-user row (with a blank spacer row
# above and below it, matching the sheet's existing layout) to the "Users"
# worksheet, and updates the active selection the way Excel leaves it after
# such an edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Users" sheet, already the active/selected tab

# Duplicate the bordered formatting of the last existing data row (17) onto
# the three rows being appended (18-20) -- a blank spacer row, the new data
# row, and another blank spacer row.
$ws.Range("A17:G17").Copy()
$ws.Range("A18:G20").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Row 19 holds the new ANZ annotation test user.
$ws.Range("A19").Value = "auAnnotationUser1"
$ws.Range("B19").Value = "Password1"
$ws.Range("E19").Value = "ANZ annotation user"
$ws.Range("G19").Value = "auannotationuser1@mailinator.com"

# Leave the selection where Excel would after typing the last value.
[void]$ws.Range("G19").Select()
